# Update cryptocurrency price (D) and volume-change (E) columns
# to match the latest scrape, per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "66.314.13"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.323.28"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.38%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "587.46"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.49%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "183.94"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.650"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +8.47%  "

$ws.Range("E9").Value = "  -2.50%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.82"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.38%  "

$ws.Range("E11").Value = "  -0.02%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "3.899.65"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("E13").Value = "  -3.25%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "66.325.90"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.77%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "26.20"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.27%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "3.329.77"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("E17").Value = "  -2.21%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "426.13"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.57%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.54"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.43%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.24"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.41"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.51%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "71.90"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.51%  "

$ws.Range("E23").Value = "  +0.12%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.69"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.463.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.515"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.54%  "

$ws.Range("E27").Value = "  +6.76%  "

$ws.Range("E28").Value = "  -3.25%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.94"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.01%  "

$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("E31").Value = "  -2.27%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "22.43"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.80%  "

$ws.Range("E33").Value = "  +0.12%  "

$ws.Range("E34").Value = "  -1.95%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "6.59"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.17%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.19"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.89%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "159.93"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("E38").Value = "  -3.22%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.885.26"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.36%  "

$ws.Range("E40").Value = "  -1.80%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "26.45"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -5.23%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.765"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("E43").Value = "  -2.48%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "40.08"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.19%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0666"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.56%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "5.93"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.78%  "

$ws.Range("E47").Value = "  -1.95%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "23.29"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -5.19%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "314.33"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.27%  "

$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("E51").Value = "  +5.61%  "
